$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(45, 8).Value = 4500  # H45: 5000 -> 4500
$ws.Cells.Item(45, 9).Value = 4000  # I45: 0 -> 4000
$ws.Cells.Item(45, 11).Value = 12000  # K45: 0 -> 12000
$ws.Cells.Item(45, 13).Value = -11808  # M45: None -> -11808
$ws.Cells.Item(63, 8).Value = 0  # H63: 69999 -> 0
$ws.Cells.Item(63, 10).Value = 0  # J63: 69999 -> 0
$ws.Cells.Item(63, 12).Value = 0  # L63: 69999 -> 0
$ws.Cells.Item(63, 14).ClearContents()  # N63: -71247 -> (deleted)
$ws.Cells.Item(66, 8).Value = 0  # H66: 69999 -> 0
$ws.Cells.Item(66, 10).Value = 0  # J66: 69999 -> 0
$ws.Cells.Item(66, 12).Value = 0  # L66: 209997 -> 0
$ws.Cells.Item(66, 14).ClearContents()  # N66: -216237 -> (deleted)
$ws.Cells.Item(70, 8).Value = 3153  # H70: 3194.182 -> 3153
$ws.Cells.Item(70, 9).Value = 3474.75  # I70: 3269.2 -> 3474.75
$ws.Cells.Item(70, 10).Value = 2969.1428  # J70: 3131.6667 -> 2969.1428
$ws.Cells.Item(70, 11).Value = 10424.25  # K70: 9807.599999999999 -> 10424.25
$ws.Cells.Item(70, 12).Value = 8907.428400000001  # L70: 9395.000100000001 -> 8907.428400000001
$ws.Cells.Item(70, 13).Value = -10154.25  # M70: -9537.599999999999 -> -10154.25
$ws.Cells.Item(70, 14).Value = -9447.428400000001  # N70: -9935.000100000001 -> -9447.428400000001
$ws.Cells.Item(73, 8).Value = 3153  # H73: 3194.182 -> 3153
$ws.Cells.Item(73, 9).Value = 3474.75  # I73: 3269.2 -> 3474.75
$ws.Cells.Item(73, 10).Value = 2969.1428  # J73: 3131.6667 -> 2969.1428
$ws.Cells.Item(73, 11).Value = 10424.25  # K73: 9807.599999999999 -> 10424.25
$ws.Cells.Item(73, 12).Value = 8907.428400000001  # L73: 9395.000100000001 -> 8907.428400000001
$ws.Cells.Item(73, 13).Value = -9488.25  # M73: -8871.599999999999 -> -9488.25
$ws.Cells.Item(73, 14).Value = -10779.4284  # N73: -11267.0001 -> -10779.4284
$ws.Cells.Item(76, 8).Value = 3987.4546  # H76: 4196.0527 -> 3987.4546
$ws.Cells.Item(76, 9).Value = 4225.4116  # I76: 4455.4 -> 4225.4116
$ws.Cells.Item(76, 10).Value = 3178.4  # J76: 3223.5 -> 3178.4
$ws.Cells.Item(76, 11).Value = 4225.4116  # K76: 4455.4 -> 4225.4116
$ws.Cells.Item(76, 12).Value = 3178.4  # L76: 3223.5 -> 3178.4
$ws.Cells.Item(76, 13).Value = -3910.4116  # M76: -4140.4 -> -3910.4116
$ws.Cells.Item(76, 14).Value = -3808.4  # N76: -3853.5 -> -3808.4
$ws.Cells.Item(79, 8).Value = 3987.4546  # H79: 4196.0527 -> 3987.4546
$ws.Cells.Item(79, 9).Value = 4225.4116  # I79: 4455.4 -> 4225.4116
$ws.Cells.Item(79, 10).Value = 3178.4  # J79: 3223.5 -> 3178.4
$ws.Cells.Item(79, 11).Value = 4225.4116  # K79: 4455.4 -> 4225.4116
$ws.Cells.Item(79, 12).Value = 3178.4  # L79: 3223.5 -> 3178.4
$ws.Cells.Item(79, 13).Value = -3133.4116  # M79: -3363.4 -> -3133.4116
$ws.Cells.Item(79, 14).Value = -5362.4  # N79: -5407.5 -> -5362.4
$ws.Cells.Item(86, 8).Value = 2693.4546  # H86: 2792.8 -> 2693.4546
$ws.Cells.Item(86, 9).Value = 2112  # I86: 2215 -> 2112
$ws.Cells.Item(86, 11).Value = 2112  # K86: 2215 -> 2112
$ws.Cells.Item(86, 13).Value = -989  # M86: -1092 -> -989
$ws.Cells.Item(89, 8).Value = 2693.4546  # H89: 2792.8 -> 2693.4546
$ws.Cells.Item(89, 9).Value = 2112  # I89: 2215 -> 2112
$ws.Cells.Item(89, 11).Value = 10560  # K89: 11075 -> 10560
$ws.Cells.Item(89, 13).Value = -4944  # M89: -5459 -> -4944
$ws.Cells.Item(112, 8).Value = 3132.3333  # H112: 3108.0908 -> 3132.3333
$ws.Cells.Item(112, 10).Value = 5749  # J112: 4832.3335 -> 5749
$ws.Cells.Item(112, 12).Value = 17247  # L112: 14497.0005 -> 17247
$ws.Cells.Item(112, 14).Value = -19463  # N112: -16713.0005 -> -19463
$ws.Cells.Item(125, 8).Value = 1498.75  # H125: 1498.3334 -> 1498.75
$ws.Cells.Item(138, 8).Value = 2748.1277  # H138: 2812.9768 -> 2748.1277
$ws.Cells.Item(138, 9).Value = 2645.2334  # I138: 2701.3667 -> 2645.2334
$ws.Cells.Item(138, 10).Value = 2929.7058  # J138: 3070.5386 -> 2929.7058
$ws.Cells.Item(138, 11).Value = 7935.7002  # K138: 8104.1001 -> 7935.7002
$ws.Cells.Item(138, 12).Value = 8789.117400000001  # L138: 9211.6158 -> 8789.117400000001
$ws.Cells.Item(138, 13).Value = -2795.7002  # M138: -2964.1001 -> -2795.7002
$ws.Cells.Item(138, 14).Value = -19069.1174  # N138: -19491.6158 -> -19069.1174
$ws.Cells.Item(141, 8).Value = 3311.4167  # H141: 3158.3635 -> 3311.4167
$ws.Cells.Item(141, 10).Value = 4996.3335  # J141: 4997 -> 4996.3335
$ws.Cells.Item(141, 12).Value = 14989.0005  # L141: 14991 -> 14989.0005
$ws.Cells.Item(141, 14).Value = -25349.0005  # N141: -25351 -> -25349.0005

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 30129.297  # H32: 30960.195 -> 30129.297
$ws.Cells.Item(32, 9).Value = 31184.854  # I32: 32123.273 -> 31184.854
$ws.Cells.Item(32, 11).Value = 31184.854  # K32: 32123.273 -> 31184.854
$ws.Cells.Item(32, 13).Value = -30897.854  # M32: -31836.273 -> -30897.854
$ws.Cells.Item(61, 8).Value = 3513.8262  # H61: 3569 -> 3513.8262
$ws.Cells.Item(61, 9).Value = 3491.7273  # I61: 3548.476 -> 3491.7273
$ws.Cells.Item(61, 11).Value = 3491.7273  # K61: 3548.476 -> 3491.7273
$ws.Cells.Item(61, 13).Value = -3279.7273  # M61: -3336.476 -> -3279.7273
$ws.Cells.Item(63, 8).Value = 5564.8  # H63: 5513.273 -> 5564.8
$ws.Cells.Item(63, 9).Value = 5081.25  # I63: 5072 -> 5081.25
$ws.Cells.Item(63, 11).Value = 5081.25  # K63: 5072 -> 5081.25
$ws.Cells.Item(63, 13).Value = -4395.25  # M63: -4386 -> -4395.25
$ws.Cells.Item(64, 8).Value = 241664.67  # H64: 243747.25 -> 241664.67
$ws.Cells.Item(64, 10).Value = 249997  # J64: 249996.33 -> 249997
$ws.Cells.Item(64, 12).Value = 249997  # L64: 249996.33 -> 249997
$ws.Cells.Item(64, 14).Value = -250493  # N64: -250492.33 -> -250493
$ws.Cells.Item(66, 8).Value = 5564.8  # H66: 5513.273 -> 5564.8
$ws.Cells.Item(66, 9).Value = 5081.25  # I66: 5072 -> 5081.25
$ws.Cells.Item(66, 11).Value = 25406.25  # K66: 25360 -> 25406.25
$ws.Cells.Item(66, 13).Value = -21974.25  # M66: -21928 -> -21974.25
$ws.Cells.Item(67, 8).Value = 241664.67  # H67: 243747.25 -> 241664.67
$ws.Cells.Item(67, 10).Value = 249997  # J67: 249996.33 -> 249997
$ws.Cells.Item(67, 12).Value = 249997  # L67: 249996.33 -> 249997
$ws.Cells.Item(67, 14).Value = -251713  # N67: -251712.33 -> -251713
$ws.Cells.Item(74, 8).Value = 2221.3076  # H74: 2396.5715 -> 2221.3076
$ws.Cells.Item(74, 9).Value = 2179.2432  # I74: 2360.0303 -> 2179.2432
$ws.Cells.Item(74, 11).Value = 2179.2432  # K74: 2360.0303 -> 2179.2432
$ws.Cells.Item(74, 13).Value = -1305.2432  # M74: -1486.0303 -> -1305.2432
$ws.Cells.Item(76, 8).Value = 37997  # H76: 37999 -> 37997
$ws.Cells.Item(76, 10).Value = 37997  # J76: 37999 -> 37997
$ws.Cells.Item(76, 12).Value = 37997  # L76: 37999 -> 37997
$ws.Cells.Item(76, 14).Value = -38673  # N76: -38675 -> -38673
$ws.Cells.Item(77, 8).Value = 2221.3076  # H77: 2396.5715 -> 2221.3076
$ws.Cells.Item(77, 9).Value = 2179.2432  # I77: 2360.0303 -> 2179.2432
$ws.Cells.Item(77, 11).Value = 10896.216  # K77: 11800.1515 -> 10896.216
$ws.Cells.Item(77, 13).Value = -6528.216  # M77: -7432.1515 -> -6528.216
$ws.Cells.Item(79, 8).Value = 37997  # H79: 37999 -> 37997
$ws.Cells.Item(79, 10).Value = 37997  # J79: 37999 -> 37997
$ws.Cells.Item(79, 12).Value = 37997  # L79: 37999 -> 37997
$ws.Cells.Item(79, 14).Value = -40337  # N79: -40339 -> -40337
$ws.Cells.Item(86, 8).Value = 22791.666  # H86: 22795 -> 22791.666
$ws.Cells.Item(86, 9).Value = 22791.666  # I86: 22795 -> 22791.666
$ws.Cells.Item(86, 11).Value = 22791.666  # K86: 22795 -> 22791.666
$ws.Cells.Item(86, 13).Value = -21605.666  # M86: -21609 -> -21605.666
$ws.Cells.Item(89, 8).Value = 22791.666  # H89: 22795 -> 22791.666
$ws.Cells.Item(89, 9).Value = 22791.666  # I89: 22795 -> 22791.666
$ws.Cells.Item(89, 11).Value = 68374.99800000001  # K89: 68385 -> 68374.99800000001
$ws.Cells.Item(89, 13).Value = -62446.99800000001  # M89: -62457 -> -62446.99800000001
$ws.Cells.Item(102, 8).Value = 2078.4138  # H102: 2121.3572 -> 2078.4138
$ws.Cells.Item(102, 9).Value = 2034.5  # I102: 2080.84 -> 2034.5
$ws.Cells.Item(102, 11).Value = 2034.5  # K102: 2080.84 -> 2034.5
$ws.Cells.Item(102, 13).Value = -412.5  # M102: -458.8400000000001 -> -412.5
$ws.Cells.Item(103, 8).Value = 99998  # H103: 100000 -> 99998
$ws.Cells.Item(103, 10).Value = 99998  # J103: 100000 -> 99998
$ws.Cells.Item(103, 12).Value = 99998  # L103: 100000 -> 99998
$ws.Cells.Item(103, 14).Value = -102342  # N103: -102344 -> -102342
$ws.Cells.Item(129, 8).Value = 73296.5  # H129: 73299.336 -> 73296.5
$ws.Cells.Item(129, 10).Value = 73296.5  # J129: 73299.336 -> 73296.5
$ws.Cells.Item(129, 12).Value = 73296.5  # L129: 73299.336 -> 73296.5
$ws.Cells.Item(129, 14).Value = -83296.5  # N129: -83299.336 -> -83296.5
$ws.Cells.Item(130, 8).Value = 0  # H130: 50429 -> 0
$ws.Cells.Item(130, 10).Value = 0  # J130: 50429 -> 0
$ws.Cells.Item(130, 12).Value = 0  # L130: 50429 -> 0
$ws.Cells.Item(130, 14).ClearContents()  # N130: -60469 -> (deleted)
$ws.Cells.Item(136, 8).Value = 3513.8262  # H136: 3569 -> 3513.8262
$ws.Cells.Item(136, 9).Value = 3491.7273  # I136: 3548.476 -> 3491.7273
$ws.Cells.Item(136, 11).Value = 10475.1819  # K136: 10645.428 -> 10475.1819
$ws.Cells.Item(136, 13).Value = -7925.1819  # M136: -8095.428 -> -7925.1819

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 3711.074  # H86: 4065.5557 -> 3711.074
$ws.Cells.Item(86, 9).Value = 3037.0908  # I86: 3524.4546 -> 3037.0908
$ws.Cells.Item(86, 10).Value = 4174.4375  # J86: 4437.5625 -> 4174.4375
$ws.Cells.Item(86, 11).Value = 3037.0908  # K86: 3524.4546 -> 3037.0908
$ws.Cells.Item(86, 12).Value = 4174.4375  # L86: 4437.5625 -> 4174.4375
$ws.Cells.Item(86, 13).Value = -1914.0908  # M86: -2401.4546 -> -1914.0908
$ws.Cells.Item(86, 14).Value = -6420.4375  # N86: -6683.5625 -> -6420.4375
$ws.Cells.Item(89, 8).Value = 3711.074  # H89: 4065.5557 -> 3711.074
$ws.Cells.Item(89, 9).Value = 3037.0908  # I89: 3524.4546 -> 3037.0908
$ws.Cells.Item(89, 10).Value = 4174.4375  # J89: 4437.5625 -> 4174.4375
$ws.Cells.Item(89, 11).Value = 15185.454  # K89: 17622.273 -> 15185.454
$ws.Cells.Item(89, 12).Value = 20872.1875  # L89: 22187.8125 -> 20872.1875
$ws.Cells.Item(89, 13).Value = -9569.454  # M89: -12006.273 -> -9569.454
$ws.Cells.Item(89, 14).Value = -32104.1875  # N89: -33419.8125 -> -32104.1875

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(6, 8).Value = 858.3158  # H6: 1499 -> 858.3158
$ws.Cells.Item(6, 9).Value = 300.8889  # I6: 1499 -> 300.8889
$ws.Cells.Item(6, 10).Value = 1360  # J6: 0 -> 1360
$ws.Cells.Item(6, 11).Value = 300.8889  # K6: 1499 -> 300.8889
$ws.Cells.Item(6, 12).Value = 1360  # L6: 0 -> 1360
$ws.Cells.Item(6, 13).Value = -187.8889  # M6: -1386 -> -187.8889
$ws.Cells.Item(6, 14).Value = -1586  # N6: None -> -1586
$ws.Cells.Item(10, 8).Value = 1403.7  # H10: 1452.9 -> 1403.7
$ws.Cells.Item(10, 9).Value = 524.6667  # I10: 710.25 -> 524.6667
$ws.Cells.Item(10, 10).Value = 2722.25  # J10: 1948 -> 2722.25
$ws.Cells.Item(10, 11).Value = 524.6667  # K10: 710.25 -> 524.6667
$ws.Cells.Item(10, 12).Value = 2722.25  # L10: 1948 -> 2722.25
$ws.Cells.Item(10, 13).Value = -385.6667  # M10: -571.25 -> -385.6667
$ws.Cells.Item(10, 14).Value = -3000.25  # N10: -2226 -> -3000.25
$ws.Cells.Item(31, 8).Value = 6831.5  # H31: 6648.4 -> 6831.5
$ws.Cells.Item(31, 9).Value = 5999.1113  # I31: 5816.636 -> 5999.1113
$ws.Cells.Item(31, 10).Value = 7663.8887  # J31: 7665 -> 7663.8887
$ws.Cells.Item(31, 11).Value = 5999.1113  # K31: 5816.636 -> 5999.1113
$ws.Cells.Item(31, 12).Value = 7663.8887  # L31: 7665 -> 7663.8887
$ws.Cells.Item(31, 13).Value = -5704.1113  # M31: -5521.636 -> -5704.1113
$ws.Cells.Item(31, 14).Value = -8253.8887  # N31: -8255 -> -8253.8887
$ws.Cells.Item(32, 8).Value = 6783.5557  # H32: 6384.1 -> 6783.5557
$ws.Cells.Item(32, 10).Value = 2867  # J32: 2847.5 -> 2867
$ws.Cells.Item(32, 12).Value = 2867  # L32: 2847.5 -> 2867
$ws.Cells.Item(32, 14).Value = -3499  # N32: -3479.5 -> -3499
$ws.Cells.Item(34, 8).Value = 6831.5  # H34: 6648.4 -> 6831.5
$ws.Cells.Item(34, 9).Value = 5999.1113  # I34: 5816.636 -> 5999.1113
$ws.Cells.Item(34, 10).Value = 7663.8887  # J34: 7665 -> 7663.8887
$ws.Cells.Item(34, 11).Value = 5999.1113  # K34: 5816.636 -> 5999.1113
$ws.Cells.Item(34, 12).Value = 7663.8887  # L34: 7665 -> 7663.8887
$ws.Cells.Item(34, 13).Value = -5797.1113  # M34: -5614.636 -> -5797.1113
$ws.Cells.Item(34, 14).Value = -8067.8887  # N34: -8069 -> -8067.8887
$ws.Cells.Item(50, 8).Value = 46974  # H50: 46647.668 -> 46974
$ws.Cells.Item(50, 10).Value = 46974  # J50: 46647.668 -> 46974
$ws.Cells.Item(50, 12).Value = 46974  # L50: 46647.668 -> 46974
$ws.Cells.Item(50, 14).Value = -48224  # N50: -47897.668 -> -48224
$ws.Cells.Item(60, 8).Value = 34799.2  # H60: 36283.715 -> 34799.2
$ws.Cells.Item(60, 10).Value = 40999  # J60: 40664.332 -> 40999
$ws.Cells.Item(60, 12).Value = 40999  # L60: 40664.332 -> 40999
$ws.Cells.Item(60, 14).Value = -42021  # N60: -41686.332 -> -42021
$ws.Cells.Item(141, 8).Value = 47662.25  # H141: 39995 -> 47662.25
$ws.Cells.Item(141, 10).Value = 55216.332  # J141: 43743.75 -> 55216.332
$ws.Cells.Item(141, 12).Value = 55216.332  # L141: 43743.75 -> 55216.332
$ws.Cells.Item(141, 14).Value = -65576.33199999999  # N141: -54103.75 -> -65576.33199999999

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(112, 8).Value = 12499.25  # H112: 8837.333000000001 -> 12499.25
$ws.Cells.Item(112, 9).Value = 9999  # I112: 6604.8 -> 9999
$ws.Cells.Item(112, 11).Value = 29997  # K112: 19814.4 -> 29997
$ws.Cells.Item(112, 13).Value = -28889  # M112: -18706.4 -> -28889
$ws.Cells.Item(131, 8).Value = 12046.2  # H131: 11541.381 -> 12046.2
$ws.Cells.Item(131, 9).Value = 865.6  # I131: 821.3333 -> 865.6
$ws.Cells.Item(131, 10).Value = 15773.066  # J131: 15829.4 -> 15773.066
$ws.Cells.Item(131, 11).Value = 2596.8  # K131: 2463.9999 -> 2596.8
$ws.Cells.Item(131, 12).Value = 47319.198  # L131: 47488.2 -> 47319.198
$ws.Cells.Item(131, 13).Value = 2443.2  # M131: 2576.0001 -> 2443.2
$ws.Cells.Item(131, 14).Value = -57399.198  # N131: -57568.2 -> -57399.198
$ws.Cells.Item(140, 8).Value = 2664.9092  # H140: 1828 -> 2664.9092
$ws.Cells.Item(140, 9).Value = 2664.9092  # I140: 1828 -> 2664.9092
$ws.Cells.Item(140, 11).Value = 7994.7276  # K140: 5484 -> 7994.7276
$ws.Cells.Item(140, 13).Value = -2814.7276  # M140: -304 -> -2814.7276

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 123.94444  # H2: 246.55 -> 123.94444
$ws.Cells.Item(2, 9).Value = 135.14285  # I2: 287 -> 135.14285
$ws.Cells.Item(2, 11).Value = 135.14285  # K2: 287 -> 135.14285
$ws.Cells.Item(2, 13).Value = -22.14285000000001  # M2: -174 -> -22.14285000000001
$ws.Cells.Item(113, 8).Value = 83396.234  # H113: 83666.62 -> 83396.234
$ws.Cells.Item(113, 9).Value = 64523.41  # I113: 64936.94 -> 64523.41
$ws.Cells.Item(113, 11).Value = 64523.41  # K113: 64936.94 -> 64523.41
$ws.Cells.Item(113, 13).Value = -62353.41  # M113: -62766.94 -> -62353.41

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 9337.25  # H7: 9675 -> 9337.25
$ws.Cells.Item(7, 9).Value = 9424.5  # I7: 9850 -> 9424.5
$ws.Cells.Item(7, 10).Value = 9250  # J7: 9500 -> 9250
$ws.Cells.Item(7, 11).Value = 9424.5  # K7: 9850 -> 9424.5
$ws.Cells.Item(7, 12).Value = 9250  # L7: 9500 -> 9250
$ws.Cells.Item(7, 13).Value = -9312.5  # M7: -9738 -> -9312.5
$ws.Cells.Item(7, 14).Value = -9474  # N7: -9724 -> -9474
$ws.Cells.Item(126, 8).Value = 9337.25  # H126: 9675 -> 9337.25
$ws.Cells.Item(126, 9).Value = 9424.5  # I126: 9850 -> 9424.5
$ws.Cells.Item(126, 10).Value = 9250  # J126: 9500 -> 9250
$ws.Cells.Item(126, 11).Value = 28273.5  # K126: 29550 -> 28273.5
$ws.Cells.Item(126, 12).Value = 27750  # L126: 28500 -> 27750
$ws.Cells.Item(126, 13).Value = -25803.5  # M126: -27080 -> -25803.5
$ws.Cells.Item(126, 14).Value = -32690  # N126: -33440 -> -32690
$ws.Cells.Item(132, 8).Value = 112913.73  # H132: 136410.44 -> 112913.73
$ws.Cells.Item(132, 9).Value = 135783.56  # I132: 152211.75 -> 135783.56
$ws.Cells.Item(132, 10).Value = 9999.5  # J132: 10000 -> 9999.5
$ws.Cells.Item(132, 11).Value = 407350.68  # K132: 456635.25 -> 407350.68
$ws.Cells.Item(132, 12).Value = 29998.5  # L132: 30000 -> 29998.5
$ws.Cells.Item(132, 13).Value = -404820.68  # M132: -454105.25 -> -404820.68
$ws.Cells.Item(132, 14).Value = -35058.5  # N132: -35060 -> -35058.5

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(18, 8).Value = 600  # H18: 110 -> 600
$ws.Cells.Item(18, 9).Value = 600  # I18: 110 -> 600
$ws.Cells.Item(18, 11).Value = 600  # K18: 110 -> 600
$ws.Cells.Item(18, 13).Value = -427  # M18: 63 -> -427
$ws.Cells.Item(132, 8).Value = 73973.42999999999  # H132: 79856.84 -> 73973.42999999999
$ws.Cells.Item(132, 9).Value = 73973.42999999999  # I132: 86362 -> 73973.42999999999
$ws.Cells.Item(132, 10).Value = 0  # J132: 1795 -> 0
$ws.Cells.Item(132, 11).Value = 221920.29  # K132: 259086 -> 221920.29
$ws.Cells.Item(132, 12).Value = 0  # L132: 5385 -> 0
$ws.Cells.Item(132, 13).Value = -219390.29  # M132: -256556 -> -219390.29
$ws.Cells.Item(132, 14).ClearContents()  # N132: -10445 -> (deleted)
